$wb = $excel.ActiveWorkbook

# --- Sheet "Overview": swap row 5 and row 6 data (A,B,E,F,G), then set new Status for row 5 ---
$wsOv = $wb.Worksheets.Item("Overview")

$a5 = $wsOv.Range("A5").Value2
$b5 = $wsOv.Range("B5").Value2
$e5 = $wsOv.Range("E5").Value2
$f5 = $wsOv.Range("F5").Value2
$g5 = $wsOv.Range("G5").Value2

$a6 = $wsOv.Range("A6").Value2
$b6 = $wsOv.Range("B6").Value2
$e6 = $wsOv.Range("E6").Value2
$f6 = $wsOv.Range("F6").Value2
$g6 = $wsOv.Range("G6").Value2

$wsOv.Range("A5").Value = $a6
$wsOv.Range("B5").Value = $b6
$wsOv.Range("E5").Value = "In Translation"
$wsOv.Range("F5").Value = "In Translation"
$wsOv.Range("G5").Value = $g6

$wsOv.Range("A6").Value = $a5
$wsOv.Range("B6").Value = $b5
$wsOv.Range("E6").Value = $e5
$wsOv.Range("F6").Value = $f5
$wsOv.Range("G6").Value = $g5

foreach ($h in $wsOv.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$B$5') {
        $h.TextToDisplay = "e2e\f533db6c-0e99-4b47-8d23-11499bbe4b5e.md"
    } elseif ($addr -eq '$B$6') {
        $h.TextToDisplay = "e2e\6ebcdc43-4f33-46d0-9a36-5d6438337773.md"
    }
}

# --- Sheet "zh-cn": swap row 5 and row 6 data (A,C,G,H), then set new Status for row 5 ---
$wsZh = $wb.Worksheets.Item("zh-cn")

$a5z = $wsZh.Range("A5").Value2
$g5z = $wsZh.Range("G5").Value2
$h5z = $wsZh.Range("H5").Value2

$a6z = $wsZh.Range("A6").Value2
$g6z = $wsZh.Range("G6").Value2
$h6z = $wsZh.Range("H6").Value2

$wsZh.Range("A5").Value = $a6z
$wsZh.Range("C5").Value = "In Translation"
$wsZh.Range("G5").Value = $g6z
$wsZh.Range("H5").Value = $h6z

$wsZh.Range("A6").Value = $a5z
$wsZh.Range("C6").Value = "Ready for handoff"
$wsZh.Range("G6").Value = $g5z
$wsZh.Range("H6").Value = $h5z

foreach ($h in $wsZh.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$A$5') {
        $h.TextToDisplay = "f533db6c-0e99-4b47-8d23-11499bbe4b5e.md"
    } elseif ($addr -eq '$A$6') {
        $h.TextToDisplay = "6ebcdc43-4f33-46d0-9a36-5d6438337773.md"
    }
}

# --- Sheet "de-de": swap row 5 and row 6 data (A,C,G,H), then set new Status for row 5 ---
$wsDe = $wb.Worksheets.Item("de-de")

$a5d = $wsDe.Range("A5").Value2
$g5d = $wsDe.Range("G5").Value2
$h5d = $wsDe.Range("H5").Value2

$a6d = $wsDe.Range("A6").Value2
$g6d = $wsDe.Range("G6").Value2
$h6d = $wsDe.Range("H6").Value2

$wsDe.Range("A5").Value = $a6d
$wsDe.Range("C5").Value = "In Translation"
$wsDe.Range("G5").Value = $g6d
$wsDe.Range("H5").Value = $h5d

$wsDe.Range("A6").Value = $a5d
$wsDe.Range("C6").Value = "Ready for handoff"
$wsDe.Range("G6").Value = $g5d
$wsDe.Range("H6").Value = $h6d

foreach ($h in $wsDe.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$A$5') {
        $h.TextToDisplay = "f533db6c-0e99-4b47-8d23-11499bbe4b5e.md"
    } elseif ($addr -eq '$A$6') {
        $h.TextToDisplay = "6ebcdc43-4f33-46d0-9a36-5d6438337773.md"
    }
}

Write-Output "done"
